$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.049662113189697
$ws.Range("B1").Value = 1.236498236656189
$ws.Range("D1").Value = 1.664909839630127
$ws.Range("E1").Value = 1.00478720664978
